$wb = $excel.ActiveWorkbook

# Sheet ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1174.3971  # H17: 1054.4625 -> 1174.3971
$ws.Cells.Item(17, 10).Value = 1174.3971  # J17: 1054.4625 -> 1174.3971
$ws.Cells.Item(17, 12).Value = 3523.1913  # L17: 3163.3875 -> 3523.1913
$ws.Cells.Item(17, 14).Value = -3859.1913  # N17: -3499.3875 -> -3859.1913

# Sheet ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 18426.285  # H18: 13159.2 -> 18426.285
$ws.Cells.Item(18, 9).Value = 5000  # I18: 801 -> 5000
$ws.Cells.Item(18, 10).Value = 20664  # J18: 16248.75 -> 20664
$ws.Cells.Item(18, 11).Value = 5000  # K18: 801 -> 5000
$ws.Cells.Item(18, 12).Value = 20664  # L18: 16248.75 -> 20664
$ws.Cells.Item(18, 13).Value = -4716  # M18: -517 -> -4716
$ws.Cells.Item(18, 14).Value = -21232  # N18: -16816.75 -> -21232

# Sheet ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 1540  # H43: 1238 -> 1540
$ws.Cells.Item(43, 9).Value = 0  # I43: 1080 -> 0
$ws.Cells.Item(43, 10).Value = 1540  # J43: 1343.3334 -> 1540
$ws.Cells.Item(43, 11).Value = 0  # K43: 1080 -> 0
$ws.Cells.Item(43, 12).Value = 1540  # L43: 1343.3334 -> 1540
$ws.Cells.Item(43, 13).ClearContents()  # M43: -1011 -> (removed)
$ws.Cells.Item(43, 14).Value = -1678  # N43: -1481.3334 -> -1678

# Sheet ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2183.5  # H98: 3293.889 -> 2183.5
$ws.Cells.Item(98, 9).Value = 2285.4707  # I98: 3293.889 -> 2285.4707
$ws.Cells.Item(98, 10).Value = 450  # J98: 0 -> 450
$ws.Cells.Item(98, 11).Value = 2285.4707  # K98: 3293.889 -> 2285.4707
$ws.Cells.Item(98, 12).Value = 450  # L98: 0 -> 450
$ws.Cells.Item(98, 13).Value = -787.4706999999999  # M98: -1795.889 -> -787.4706999999999
$ws.Cells.Item(98, 14).Value = -3446  # N98: None -> -3446

# Sheet ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 5229.1055  # H112: 4854.857 -> 5229.1055
$ws.Cells.Item(112, 10).Value = 5661.9414  # J112: 5202.737 -> 5661.9414
$ws.Cells.Item(112, 12).Value = 16985.8242  # L112: 15608.211 -> 16985.8242
$ws.Cells.Item(112, 14).Value = -19201.8242  # N112: -17824.211 -> -19201.8242

# Sheet ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 12794.9  # H113: 11867.272 -> 12794.9
$ws.Cells.Item(113, 9).Value = 13994.333  # I113: 12854 -> 13994.333
$ws.Cells.Item(113, 11).Value = 13994.333  # K113: 12854 -> 13994.333
$ws.Cells.Item(113, 13).Value = -10740.333  # M113: -9600 -> -10740.333

# Sheet ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 2183.5  # H122: 3293.889 -> 2183.5
$ws.Cells.Item(122, 9).Value = 2285.4707  # I122: 3293.889 -> 2285.4707
$ws.Cells.Item(122, 10).Value = 450  # J122: 0 -> 450
$ws.Cells.Item(122, 11).Value = 6856.4121  # K122: 9881.667000000001 -> 6856.4121
$ws.Cells.Item(122, 12).Value = 1350  # L122: 0 -> 1350
$ws.Cells.Item(122, 13).Value = -4406.4121  # M122: -7431.667000000001 -> -4406.4121
$ws.Cells.Item(122, 14).Value = -6250  # N122: None -> -6250

# Sheet ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1098.3478  # H132: 1119.4222 -> 1098.3478
$ws.Cells.Item(132, 9).Value = 931.4048  # I132: 950.46344 -> 931.4048
$ws.Cells.Item(132, 11).Value = 2794.2144  # K132: 2851.39032 -> 2794.2144
$ws.Cells.Item(132, 13).Value = -264.2143999999998  # M132: -321.39032 -> -264.2143999999998

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1564.9524  # H137: 1583.3334 -> 1564.9524
$ws.Cells.Item(137, 9).Value = 1456.1428  # I137: 1476.5714 -> 1456.1428
$ws.Cells.Item(137, 10).Value = 1782.5714  # J137: 1796.8572 -> 1782.5714
$ws.Cells.Item(137, 11).Value = 4368.428400000001  # K137: 4429.7142 -> 4368.428400000001
$ws.Cells.Item(137, 12).Value = 5347.7142  # L137: 5390.571599999999 -> 5347.7142
$ws.Cells.Item(137, 13).Value = -1818.428400000001  # M137: -1879.7142 -> -1818.428400000001
$ws.Cells.Item(137, 14).Value = -10447.7142  # N137: -10490.5716 -> -10447.7142

# Sheet ALC row 140
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(140, 8).Value = 48817.855  # H140: 49325.855 -> 48817.855
$ws.Cells.Item(140, 10).Value = 48817.855  # J140: 49325.855 -> 48817.855
$ws.Cells.Item(140, 12).Value = 48817.855  # L140: 49325.855 -> 48817.855
$ws.Cells.Item(140, 14).Value = -59177.855  # N140: -59685.855 -> -59177.855

# Sheet ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2585023  # H2: 3877118.8 -> 2585023
$ws.Cells.Item(2, 9).Value = 3323172.5  # I2: 7752804 -> 3323172.5
$ws.Cells.Item(2, 10).Value = 1500  # J2: 1433.3334 -> 1500
$ws.Cells.Item(2, 11).Value = 3323172.5  # K2: 7752804 -> 3323172.5
$ws.Cells.Item(2, 12).Value = 1500  # L2: 1433.3334 -> 1500
$ws.Cells.Item(2, 13).Value = -3323059.5  # M2: -7752691 -> -3323059.5
$ws.Cells.Item(2, 14).Value = -1726  # N2: -1659.3334 -> -1726

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2888.2207  # H32: 3016.6177 -> 2888.2207
$ws.Cells.Item(32, 9).Value = 2018.9464  # I32: 2185.491 -> 2018.9464
$ws.Cells.Item(32, 10).Value = 6944.8335  # J32: 6532.923 -> 6944.8335
$ws.Cells.Item(32, 11).Value = 2018.9464  # K32: 2185.491 -> 2018.9464
$ws.Cells.Item(32, 12).Value = 6944.8335  # L32: 6532.923 -> 6944.8335
$ws.Cells.Item(32, 13).Value = -1731.9464  # M32: -1898.491 -> -1731.9464
$ws.Cells.Item(32, 14).Value = -7518.8335  # N32: -7106.923 -> -7518.8335

# Sheet ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1329.3429  # H74: 1341.5143 -> 1329.3429
$ws.Cells.Item(74, 9).Value = 1097.1034  # I74: 1111.7931 -> 1097.1034
$ws.Cells.Item(74, 11).Value = 1097.1034  # K74: 1111.7931 -> 1097.1034
$ws.Cells.Item(74, 13).Value = -223.1034  # M74: -237.7931000000001 -> -223.1034

# Sheet ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1329.3429  # H77: 1341.5143 -> 1329.3429
$ws.Cells.Item(77, 9).Value = 1097.1034  # I77: 1111.7931 -> 1097.1034
$ws.Cells.Item(77, 11).Value = 5485.517  # K77: 5558.9655 -> 5485.517
$ws.Cells.Item(77, 13).Value = -1117.517  # M77: -1190.9655 -> -1117.517

# Sheet ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 2063.647  # H102: 2133.1333 -> 2063.647
$ws.Cells.Item(102, 9).Value = 1775.6154  # I102: 1863.4546 -> 1775.6154
$ws.Cells.Item(102, 10).Value = 2999.75  # J102: 2874.75 -> 2999.75
$ws.Cells.Item(102, 11).Value = 1775.6154  # K102: 1863.4546 -> 1775.6154
$ws.Cells.Item(102, 12).Value = 2999.75  # L102: 2874.75 -> 2999.75
$ws.Cells.Item(102, 13).Value = -153.6153999999999  # M102: -241.4546 -> -153.6153999999999
$ws.Cells.Item(102, 14).Value = -6243.75  # N102: -6118.75 -> -6243.75

# Sheet ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 2585023  # H116: 3877118.8 -> 2585023
$ws.Cells.Item(116, 9).Value = 3323172.5  # I116: 7752804 -> 3323172.5
$ws.Cells.Item(116, 10).Value = 1500  # J116: 1433.3334 -> 1500
$ws.Cells.Item(116, 11).Value = 3323172.5  # K116: 7752804 -> 3323172.5
$ws.Cells.Item(116, 12).Value = 1500  # L116: 1433.3334 -> 1500
$ws.Cells.Item(116, 13).Value = -3320878.5  # M116: -7750510 -> -3320878.5
$ws.Cells.Item(116, 14).Value = -6088  # N116: -6021.3334 -> -6088

# Sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 10000  # H122: 2858.375 -> 10000
$ws.Cells.Item(122, 9).Value = 0  # I122: 1838.1428 -> 0
$ws.Cells.Item(122, 11).Value = 0  # K122: 5514.428400000001 -> 0
$ws.Cells.Item(122, 13).ClearContents()  # M122: -3064.428400000001 -> (removed)

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1566.4117  # H132: 1637.4062 -> 1566.4117
$ws.Cells.Item(132, 9).Value = 1023.75  # I132: 1096.1428 -> 1023.75
$ws.Cells.Item(132, 10).Value = 2868.8  # J132: 2670.7273 -> 2868.8
$ws.Cells.Item(132, 11).Value = 3071.25  # K132: 3288.4284 -> 3071.25
$ws.Cells.Item(132, 12).Value = 8606.400000000001  # L132: 8012.1819 -> 8606.400000000001
$ws.Cells.Item(132, 13).Value = -541.25  # M132: -758.4284000000002 -> -541.25
$ws.Cells.Item(132, 14).Value = -13666.4  # N132: -13072.1819 -> -13666.4

# Sheet BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2585023  # H3: 3877118.8 -> 2585023
$ws.Cells.Item(3, 9).Value = 3323172.5  # I3: 7752804 -> 3323172.5
$ws.Cells.Item(3, 10).Value = 1500  # J3: 1433.3334 -> 1500
$ws.Cells.Item(3, 11).Value = 3323172.5  # K3: 7752804 -> 3323172.5
$ws.Cells.Item(3, 12).Value = 1500  # L3: 1433.3334 -> 1500
$ws.Cells.Item(3, 13).Value = -3323058.5  # M3: -7752690 -> -3323058.5
$ws.Cells.Item(3, 14).Value = -1728  # N3: -1661.3334 -> -1728

# Sheet BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 98081.62  # H86: 79520.30499999999 -> 98081.62
$ws.Cells.Item(86, 9).Value = 2290.3  # I86: 2047.8 -> 2290.3
$ws.Cells.Item(86, 11).Value = 2290.3  # K86: 2047.8 -> 2290.3
$ws.Cells.Item(86, 13).Value = -1167.3  # M86: -924.8 -> -1167.3

# Sheet BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 98081.62  # H89: 79520.30499999999 -> 98081.62
$ws.Cells.Item(89, 9).Value = 2290.3  # I89: 2047.8 -> 2290.3
$ws.Cells.Item(89, 11).Value = 11451.5  # K89: 10239 -> 11451.5
$ws.Cells.Item(89, 13).Value = -5835.5  # M89: -4623 -> -5835.5

# Sheet BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 656.4167  # H94: 670.2083 -> 656.4167
$ws.Cells.Item(94, 9).Value = 473.0476  # I94: 488.8095 -> 473.0476
$ws.Cells.Item(94, 11).Value = 473.0476  # K94: 488.8095 -> 473.0476
$ws.Cells.Item(94, 13).Value = -22.04759999999999  # M94: -37.80950000000001 -> -22.04759999999999

# Sheet CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 812.5  # H16: 836.25 -> 812.5
$ws.Cells.Item(16, 9).Value = 772.4  # I16: 777.3333 -> 772.4
$ws.Cells.Item(16, 11).Value = 772.4  # K16: 777.3333 -> 772.4
$ws.Cells.Item(16, 13).Value = -485.4  # M16: -490.3333 -> -485.4

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1623.5555  # H31: 1816 -> 1623.5555
$ws.Cells.Item(31, 9).Value = 999.5  # I31: 999.3333 -> 999.5
$ws.Cells.Item(31, 10).Value = 2122.8  # J31: 2428.5 -> 2122.8
$ws.Cells.Item(31, 11).Value = 999.5  # K31: 999.3333 -> 999.5
$ws.Cells.Item(31, 12).Value = 2122.8  # L31: 2428.5 -> 2122.8
$ws.Cells.Item(31, 13).Value = -704.5  # M31: -704.3333 -> -704.5
$ws.Cells.Item(31, 14).Value = -2712.8  # N31: -3018.5 -> -2712.8

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1623.5555  # H34: 1816 -> 1623.5555
$ws.Cells.Item(34, 9).Value = 999.5  # I34: 999.3333 -> 999.5
$ws.Cells.Item(34, 10).Value = 2122.8  # J34: 2428.5 -> 2122.8
$ws.Cells.Item(34, 11).Value = 999.5  # K34: 999.3333 -> 999.5
$ws.Cells.Item(34, 12).Value = 2122.8  # L34: 2428.5 -> 2122.8
$ws.Cells.Item(34, 13).Value = -797.5  # M34: -797.3333 -> -797.5
$ws.Cells.Item(34, 14).Value = -2526.8  # N34: -2832.5 -> -2526.8

# Sheet CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 424.26315  # H107: 384.94446 -> 424.26315
$ws.Cells.Item(107, 9).Value = 321.5  # I107: 311.26666 -> 321.5
$ws.Cells.Item(107, 10).Value = 712  # J107: 753.3333 -> 712
$ws.Cells.Item(107, 11).Value = 321.5  # K107: 311.26666 -> 321.5
$ws.Cells.Item(107, 12).Value = 712  # L107: 753.3333 -> 712
$ws.Cells.Item(107, 13).Value = 1598.5  # M107: 1608.73334 -> 1598.5
$ws.Cells.Item(107, 14).Value = -4552  # N107: -4593.3333 -> -4552

# Sheet CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 812.5  # H113: 836.25 -> 812.5
$ws.Cells.Item(113, 9).Value = 772.4  # I113: 777.3333 -> 772.4
$ws.Cells.Item(113, 11).Value = 772.4  # K113: 777.3333 -> 772.4
$ws.Cells.Item(113, 13).Value = 1397.6  # M113: 1392.6667 -> 1397.6

# Sheet CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 9158.416999999999  # H113: 9907.362999999999 -> 9158.416999999999
$ws.Cells.Item(113, 9).Value = 17260.166  # I113: 17265.166 -> 17260.166
$ws.Cells.Item(113, 10).Value = 1056.6666  # J113: 1078 -> 1056.6666
$ws.Cells.Item(113, 11).Value = 51780.49800000001  # K113: 51795.49800000001 -> 51780.49800000001
$ws.Cells.Item(113, 12).Value = 3169.9998  # L113: 3234 -> 3169.9998
$ws.Cells.Item(113, 13).Value = -49610.49800000001  # M113: -49625.49800000001 -> -49610.49800000001
$ws.Cells.Item(113, 14).Value = -7509.9998  # N113: -7574 -> -7509.9998

# Sheet CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 793.17  # H131: 7834.758 -> 793.17
$ws.Cells.Item(131, 9).Value = 377.22223  # I131: 472.7 -> 377.22223
$ws.Cells.Item(131, 10).Value = 834.3077  # J131: 8661.955 -> 834.3077
$ws.Cells.Item(131, 11).Value = 1131.66669  # K131: 1418.1 -> 1131.66669
$ws.Cells.Item(131, 12).Value = 2502.9231  # L131: 25985.865 -> 2502.9231
$ws.Cells.Item(131, 13).Value = 3908.33331  # M131: 3621.9 -> 3908.33331
$ws.Cells.Item(131, 14).Value = -12582.9231  # N131: -36065.865 -> -12582.9231

# Sheet CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 3147.1667  # H139: 2228.1 -> 3147.1667
$ws.Cells.Item(139, 9).Value = 4940  # I139: 2255.4 -> 4940
$ws.Cells.Item(139, 10).Value = 2250.75  # J139: 2200.8 -> 2250.75
$ws.Cells.Item(139, 11).Value = 14820  # K139: 6766.200000000001 -> 14820
$ws.Cells.Item(139, 12).Value = 6752.25  # L139: 6602.400000000001 -> 6752.25
$ws.Cells.Item(139, 13).Value = -9680  # M139: -1626.200000000001 -> -9680
$ws.Cells.Item(139, 14).Value = -17032.25  # N139: -16882.4 -> -17032.25

# Sheet GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3780.375  # H102: 3475.3 -> 3780.375
$ws.Cells.Item(102, 9).Value = 3891.8572  # I102: 3656.625 -> 3891.8572
$ws.Cells.Item(102, 10).Value = 3000  # J102: 2750 -> 3000
$ws.Cells.Item(102, 11).Value = 3891.8572  # K102: 3656.625 -> 3891.8572
$ws.Cells.Item(102, 12).Value = 3000  # L102: 2750 -> 3000
$ws.Cells.Item(102, 13).Value = -2269.8572  # M102: -2034.625 -> -2269.8572
$ws.Cells.Item(102, 14).Value = -6244  # N102: -5994 -> -6244

# Sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1473.4  # H122: 1387.4062 -> 1473.4
$ws.Cells.Item(122, 9).Value = 1145.375  # I122: 1125 -> 1145.375
$ws.Cells.Item(122, 10).Value = 1848.2858  # J122: 1684.8 -> 1848.2858
$ws.Cells.Item(122, 11).Value = 3436.125  # K122: 3375 -> 3436.125
$ws.Cells.Item(122, 12).Value = 5544.857400000001  # L122: 5054.4 -> 5544.857400000001
$ws.Cells.Item(122, 13).Value = -986.125  # M122: -925 -> -986.125
$ws.Cells.Item(122, 14).Value = -10444.8574  # N122: -9954.4 -> -10444.8574

# Sheet GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2461590.5  # H126: 1826748.8 -> 2461590.5
$ws.Cells.Item(126, 9).Value = 2780583.5  # I126: 1986582.1 -> 2780583.5
$ws.Cells.Item(126, 11).Value = 8341750.5  # K126: 5959746.300000001 -> 8341750.5
$ws.Cells.Item(126, 13).Value = -8339280.5  # M126: -5957276.300000001 -> -8339280.5

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 896548.0600000001  # H132: 1167988.8 -> 896548.0600000001
$ws.Cells.Item(132, 9).Value = 1375320.2  # I132: 1750215.6 -> 1375320.2
$ws.Cells.Item(132, 10).Value = 2839.9333  # J132: 3535.0908 -> 2839.9333
$ws.Cells.Item(132, 11).Value = 4125960.6  # K132: 5250646.800000001 -> 4125960.6
$ws.Cells.Item(132, 12).Value = 8519.7999  # L132: 10605.2724 -> 8519.7999
$ws.Cells.Item(132, 13).Value = -4123430.6  # M132: -5248116.800000001 -> -4123430.6
$ws.Cells.Item(132, 14).Value = -13579.7999  # N132: -15665.2724 -> -13579.7999

# Sheet GSM row 139
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(139, 8).Value = 50945.715  # H139: 57545 -> 50945.715
$ws.Cells.Item(139, 10).Value = 50945.715  # J139: 57545 -> 50945.715
$ws.Cells.Item(139, 12).Value = 50945.715  # L139: 57545 -> 50945.715
$ws.Cells.Item(139, 14).Value = -61225.715  # N139: -67825 -> -61225.715

# Sheet LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3028.5715  # H22: 1926.3636 -> 3028.5715
$ws.Cells.Item(22, 9).Value = 5300  # I22: 1961.25 -> 5300
$ws.Cells.Item(22, 10).Value = 2120  # J22: 1833.3334 -> 2120
$ws.Cells.Item(22, 11).Value = 5300  # K22: 1961.25 -> 5300
$ws.Cells.Item(22, 12).Value = 2120  # L22: 1833.3334 -> 2120
$ws.Cells.Item(22, 13).Value = -5005  # M22: -1666.25 -> -5005
$ws.Cells.Item(22, 14).Value = -2710  # N22: -2423.3334 -> -2710

# Sheet LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 3028.5715  # H27: 1926.3636 -> 3028.5715
$ws.Cells.Item(27, 9).Value = 5300  # I27: 1961.25 -> 5300
$ws.Cells.Item(27, 10).Value = 2120  # J27: 1833.3334 -> 2120
$ws.Cells.Item(27, 11).Value = 5300  # K27: 1961.25 -> 5300
$ws.Cells.Item(27, 12).Value = 2120  # L27: 1833.3334 -> 2120
$ws.Cells.Item(27, 13).Value = -5193  # M27: -1854.25 -> -5193
$ws.Cells.Item(27, 14).Value = -2334  # N27: -2047.3334 -> -2334

# Sheet LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2953.4285  # H46: 2930.5 -> 2953.4285
$ws.Cells.Item(46, 9).Value = 1500  # I46: 1400 -> 1500
$ws.Cells.Item(46, 10).Value = 3195.6667  # J46: 3149.1428 -> 3195.6667
$ws.Cells.Item(46, 11).Value = 1500  # K46: 1400 -> 1500
$ws.Cells.Item(46, 12).Value = 3195.6667  # L46: 3149.1428 -> 3195.6667
$ws.Cells.Item(46, 13).Value = -1312  # M46: -1212 -> -1312
$ws.Cells.Item(46, 14).Value = -3571.6667  # N46: -3525.1428 -> -3571.6667

# Sheet LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3024.625  # H61: 3074.625 -> 3024.625
$ws.Cells.Item(61, 9).Value = 2742.4285  # I61: 3016.1667 -> 2742.4285
$ws.Cells.Item(61, 10).Value = 5000  # J61: 3250 -> 5000
$ws.Cells.Item(61, 11).Value = 2742.4285  # K61: 3016.1667 -> 2742.4285
$ws.Cells.Item(61, 12).Value = 5000  # L61: 3250 -> 5000
$ws.Cells.Item(61, 13).Value = -2540.4285  # M61: -2814.1667 -> -2540.4285
$ws.Cells.Item(61, 14).Value = -5404  # N61: -3654 -> -5404

# Sheet LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 3267.375  # H68: 3534.1428 -> 3267.375
$ws.Cells.Item(68, 9).Value = 3019.8572  # I68: 3289.8333 -> 3019.8572
$ws.Cells.Item(68, 11).Value = 3019.8572  # K68: 3289.8333 -> 3019.8572
$ws.Cells.Item(68, 13).Value = -2270.8572  # M68: -2540.8333 -> -2270.8572

# Sheet LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 3267.375  # H71: 3534.1428 -> 3267.375
$ws.Cells.Item(71, 9).Value = 3019.8572  # I71: 3289.8333 -> 3019.8572
$ws.Cells.Item(71, 11).Value = 15099.286  # K71: 16449.1665 -> 15099.286
$ws.Cells.Item(71, 13).Value = -11355.286  # M71: -12705.1665 -> -11355.286

# Sheet LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 3024.625  # H113: 3074.625 -> 3024.625
$ws.Cells.Item(113, 9).Value = 2742.4285  # I113: 3016.1667 -> 2742.4285
$ws.Cells.Item(113, 10).Value = 5000  # J113: 3250 -> 5000
$ws.Cells.Item(113, 11).Value = 2742.4285  # K113: 3016.1667 -> 2742.4285
$ws.Cells.Item(113, 12).Value = 5000  # L113: 3250 -> 5000
$ws.Cells.Item(113, 13).Value = -572.4285  # M113: -846.1667000000002 -> -572.4285
$ws.Cells.Item(113, 14).Value = -9340  # N113: -7590 -> -9340

# Sheet WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1626.6  # H81: 1672.5834 -> 1626.6
$ws.Cells.Item(81, 9).Value = 1626.6  # I81: 1696.7778 -> 1626.6
$ws.Cells.Item(81, 10).Value = 0  # J81: 1600 -> 0
$ws.Cells.Item(81, 11).Value = 3253.2  # K81: 3393.5556 -> 3253.2
$ws.Cells.Item(81, 12).Value = 0  # L81: 3200 -> 0
$ws.Cells.Item(81, 13).Value = -2192.2  # M81: -2332.5556 -> -2192.2
$ws.Cells.Item(81, 14).ClearContents()  # N81: -5322 -> (removed)

# Sheet WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 1626.6  # H84: 1672.5834 -> 1626.6
$ws.Cells.Item(84, 9).Value = 1626.6  # I84: 1696.7778 -> 1626.6
$ws.Cells.Item(84, 10).Value = 0  # J84: 1600 -> 0
$ws.Cells.Item(84, 11).Value = 16266  # K84: 16967.778 -> 16266
$ws.Cells.Item(84, 12).Value = 0  # L84: 16000 -> 0
$ws.Cells.Item(84, 13).Value = -10962  # M84: -11663.778 -> -10962
$ws.Cells.Item(84, 14).ClearContents()  # N84: -26608 -> (removed)

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1122  # H132: 1142.0488 -> 1122
$ws.Cells.Item(132, 9).Value = 816.3125  # I132: 832.9677 -> 816.3125
$ws.Cells.Item(132, 11).Value = 2448.9375  # K132: 2498.9031 -> 2448.9375
$ws.Cells.Item(132, 13).Value = 81.0625  # M132: 31.09690000000001 -> 81.0625

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 15875370  # H136: 18521192 -> 15875370
$ws.Cells.Item(136, 9).Value = 19843502  # I136: 25255188 -> 19843502
$ws.Cells.Item(136, 10).Value = 2841.4285  # J136: 2698.75 -> 2841.4285
$ws.Cells.Item(136, 11).Value = 59530506  # K136: 75765564 -> 59530506
$ws.Cells.Item(136, 12).Value = 8524.2855  # L136: 8096.25 -> 8524.2855
$ws.Cells.Item(136, 13).Value = -59527956  # M136: -75763014 -> -59527956
$ws.Cells.Item(136, 14).Value = -13624.2855  # N136: -13196.25 -> -13624.2855
